$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 37448
$ws.Range("J3").Value = 37448
$ws.Range("L3").Value = 37448
$ws.Range("N3").Value = -37676

$ws.Range("H68").Value = 50295
$ws.Range("J68").Value = 50295
$ws.Range("L68").Value = 50295
$ws.Range("N68").Value = -51793

$ws.Range("H71").Value = 50295
$ws.Range("J71").Value = 50295
$ws.Range("L71").Value = 150885
$ws.Range("N71").Value = -158373

$ws.Range("H98").Value = 689.8889
$ws.Range("I98").Value = 689.8889
$ws.Range("K98").Value = 689.8889
$ws.Range("M98").Value = 808.1111

$ws.Range("H102").Value = 37448
$ws.Range("J102").Value = 37448
$ws.Range("L102").Value = 37448
$ws.Range("N102").Value = -43938

$ws.Range("H122").Value = 689.8889
$ws.Range("I122").Value = 689.8889
$ws.Range("K122").Value = 2069.6667
$ws.Range("M122").Value = 380.3332999999998

$ws.Range("H141").Value = 8750
$ws.Range("I141").Value = 7500
$ws.Range("K141").Value = 22500
$ws.Range("M141").Value = -17320

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 12860428
$ws.Range("I6").Value = 12004200
$ws.Range("K6").Value = 12004200
$ws.Range("M6").Value = -12004027

$ws.Range("H32").Value = 33335344
$ws.Range("I32").Value = 3015
$ws.Range("K32").Value = 3015
$ws.Range("M32").Value = -2728

$ws.Range("H74").Value = 2723.2856
$ws.Range("I74").Value = 2308.7
$ws.Range("K74").Value = 2308.7
$ws.Range("M74").Value = -1434.7

$ws.Range("H77").Value = 2723.2856
$ws.Range("I77").Value = 2308.7
$ws.Range("K77").Value = 11543.5
$ws.Range("M77").Value = -7175.5

$ws.Range("H97").Value = 741.8461
$ws.Range("J97").Value = 410
$ws.Range("L97").Value = 410
$ws.Range("N97").Value = -1402

$ws.Range("H102").Value = 6254755
$ws.Range("I102").Value = 8336639.5
$ws.Range("K102").Value = 8336639.5
$ws.Range("M102").Value = -8335017.5

$ws.Range("H132").Value = 1792.7778
$ws.Range("I132").Value = 1792.7778
$ws.Range("K132").Value = 5378.3334
$ws.Range("M132").Value = -2848.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 8000
$ws.Range("I31").Value = 8000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 8000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -7748
$ws.Range("N31").ClearContents()

$ws.Range("H134").Value = 5764.2
$ws.Range("J134").Value = 30400
$ws.Range("L134").Value = 91200
$ws.Range("N134").Value = -96270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 843
$ws.Range("I22").Value = 843
$ws.Range("K22").Value = 843
$ws.Range("M22").Value = -493

$ws.Range("H31").Value = 5632.3145
$ws.Range("I31").Value = 3507.2778
$ws.Range("K31").Value = 3507.2778
$ws.Range("M31").Value = -3212.2778

$ws.Range("H34").Value = 5632.3145
$ws.Range("I34").Value = 3507.2778
$ws.Range("K34").Value = 3507.2778
$ws.Range("M34").Value = -3305.2778

$ws.Range("H43").Value = 19696
$ws.Range("J43").Value = 19696
$ws.Range("L43").Value = 19696
$ws.Range("N43").Value = -20064

$ws.Range("H101").Value = 19696
$ws.Range("J101").Value = 19696
$ws.Range("L101").Value = 19696
$ws.Range("N101").Value = -26186

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2307.6924
$ws.Range("I5").Value = 1628.5714
$ws.Range("K5").Value = 4885.7142
$ws.Range("M5").Value = -4773.7142

$ws.Range("H12").Value = 178.58824
$ws.Range("I12").Value = 167.85715
$ws.Range("J12").Value = 186.1
$ws.Range("K12").Value = 503.57145
$ws.Range("L12").Value = 558.3
$ws.Range("M12").Value = -330.57145
$ws.Range("N12").Value = -904.3

$ws.Range("H113").Value = 1687.5264
$ws.Range("I113").Value = 1025.4286
$ws.Range("J113").Value = 2073.75
$ws.Range("K113").Value = 3076.2858
$ws.Range("L113").Value = 6221.25
$ws.Range("M113").Value = -906.2857999999997
$ws.Range("N113").Value = -10561.25

$ws.Range("H135").Value = 2307.6924
$ws.Range("I135").Value = 1628.5714
$ws.Range("K135").Value = 14657.1426
$ws.Range("M135").Value = -12122.1426

$ws.Range("H137").Value = 3479.4
$ws.Range("J137").Value = 5249.6665
$ws.Range("L137").Value = 15748.9995
$ws.Range("N137").Value = -25948.9995

$ws.Range("H140").Value = 2271.8462
$ws.Range("I140").Value = 1828.6666
$ws.Range("J140").Value = 4133.2
$ws.Range("K140").Value = 5485.9998
$ws.Range("L140").Value = 12399.6
$ws.Range("M140").Value = -305.9997999999996
$ws.Range("N140").Value = -22759.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 732.3333
$ws.Range("J14").Value = 732.3333
$ws.Range("L14").Value = 732.3333
$ws.Range("N14").Value = -1068.3333

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6550
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1601.4286
$ws.Range("I22").Value = 927.5
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 927.5
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -632.5
$ws.Range("N22").Value = -3090

$ws.Range("H27").Value = 1601.4286
$ws.Range("I27").Value = 927.5
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 927.5
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -820.5
$ws.Range("N27").Value = -2714

$ws.Range("H46").Value = 6204.0835
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 6204.0835
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 6204.0835
$ws.Range("N46").Value = -6580.0835
$ws.Range("M46").ClearContents()

$ws.Range("H122").Value = 3015.1538
$ws.Range("J122").Value = 2999.5
$ws.Range("L122").Value = 8998.5
$ws.Range("N122").Value = -13898.5

$ws.Range("H132").Value = 5861.857
$ws.Range("I132").Value = 6377.1577
$ws.Range("J132").Value = 966.5
$ws.Range("K132").Value = 19131.4731
$ws.Range("L132").Value = 2899.5
$ws.Range("M132").Value = -16601.4731
$ws.Range("N132").Value = -7959.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992

$ws.Range("H96").Value = 1811.8889
$ws.Range("I96").Value = 1786.7142
$ws.Range("K96").Value = 1786.7142
$ws.Range("M96").Value = -413.7141999999999

$ws.Range("H126").Value = 4158.5264
$ws.Range("I126").Value = 2072.4
$ws.Range("K126").Value = 6217.200000000001
$ws.Range("M126").Value = -3747.200000000001

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
